# This script applies the weekly update to the "Ajo" (garlic) sheet.
# Effect: four brand-new records (two market days) are inserted right
# after row 50. Since the runtime doesn't insert physical rows, we
# shift the existing data block (rows 51-96) down by four rows (to
# 55-100) and then write the four new records into rows 51-54. The
# sheet's used range / dimension grows from A1:R96 to A1:R100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 18   # column R
$oldFirst = 51
$oldLast = 96
$shift = 4

# Walk bottom-up so we never overwrite a source row before reading it.
for ($r = $oldLast; $r -ge $oldFirst; $r--) {
    $destRow = $r + $shift
    for ($c = 1; $c -le $lastCol; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $val = $srcCell.Value2
        $destCell = $ws.Cells.Item($destRow, $c)
        $destCell.Value = $val
        if ($c -eq 4) {
            # Column D holds dates; reapply the date number format since
            # newly-touched cells otherwise fall back to the default style.
            $destCell.NumberFormat = $srcCell.NumberFormat
        }
    }
}

# Helper to write one full data record (columns A-R) into a given row.
function Set-AjoRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($Row, 1).Value = 9
    $ws.Cells.Item($Row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($Row, 3).Value = "Metropolitana"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 13
    $ws.Cells.Item($Row, 6).Value = 100112003
    $ws.Cells.Item($Row, 7).Value = "Ajo"
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Four new records for market dates 2021-08-27 (serial 44435) and
# 2021-08-23 (serial 44431), Chino / Primera quality.
Set-AjoRow 51 44435 "Chino" "Primera" 430 15000 15500 15250 "$/caja 10 kilos" "China" 1525 10
Set-AjoRow 52 44435 "Chino" "Primera" 340 14500 15000 14750 "$/malla 10 kilos" "China" 1475 10
Set-AjoRow 53 44431 "Chino" "Primera" 430 15000 15500 15250 "$/caja 10 kilos" "China" 1525 10
Set-AjoRow 54 44431 "Chino" "Primera" 340 14500 15000 14750 "$/malla 10 kilos" "China" 1475 10
